$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 7704
$ws.Cells.Item(3, 6).Value = 7704
$ws.Cells.Item(5, 6).Value = 7881
$ws.Cells.Item(6, 6).Value = 40
$ws.Cells.Item(8, 6).Value = 32
$ws.Cells.Item(9, 6).Value = 6657
$ws.Cells.Item(10, 6).Value = 3384
$ws.Cells.Item(12, 6).Value = 3724
$ws.Cells.Item(14, 6).Value = 47
$ws.Cells.Item(17, 6).Value = 70
$ws.Cells.Item(20, 6).Value = 36
$ws.Cells.Item(21, 6).Value = 314
$ws.Cells.Item(22, 6).Value = 6
$ws.Cells.Item(24, 6).Value = 3847
$ws.Cells.Item(25, 6).Value = 117
$ws.Cells.Item(27, 6).Value = 956
$ws.Cells.Item(29, 6).Value = 1484
$ws.Cells.Item(31, 6).Value = 54
$ws.Cells.Item(32, 6).Value = 2757
$ws.Cells.Item(33, 6).Value = 1833
$ws.Cells.Item(35, 6).Value = 48
$ws.Cells.Item(36, 6).Value = 59
$ws.Cells.Item(37, 6).Value = 3690
$ws.Cells.Item(38, 6).Value = 319
$ws.Cells.Item(40, 6).Value = 44
$ws.Cells.Item(41, 6).Value = 918
$ws.Cells.Item(42, 6).Value = 539
$ws.Cells.Item(44, 6).Value = 1425
$ws.Cells.Item(45, 6).Value = 244
$ws.Cells.Item(48, 6).Value = 642

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 240
$ws.Cells.Item(6, 6).Value = 408
$ws.Cells.Item(13, 6).Value = 90
$ws.Cells.Item(16, 6).Value = 16
$ws.Cells.Item(17, 6).Value = 34

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 135

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 135
$ws.Cells.Item(3, 6).Value = 240
$ws.Cells.Item(5, 6).Value = 7704
$ws.Cells.Item(6, 6).Value = 7704
$ws.Cells.Item(7, 6).Value = 7881
$ws.Cells.Item(8, 6).Value = 40
$ws.Cells.Item(9, 6).Value = 32
$ws.Cells.Item(10, 6).Value = 6657
$ws.Cells.Item(11, 6).Value = 3384
$ws.Cells.Item(12, 6).Value = 3724
$ws.Cells.Item(13, 6).Value = 47
$ws.Cells.Item(16, 6).Value = 70
$ws.Cells.Item(19, 6).Value = 36
$ws.Cells.Item(20, 6).Value = 314
$ws.Cells.Item(23, 6).Value = 3847
$ws.Cells.Item(25, 6).Value = 117
$ws.Cells.Item(28, 6).Value = 956
$ws.Cells.Item(30, 6).Value = 1484
$ws.Cells.Item(32, 6).Value = 54
$ws.Cells.Item(33, 6).Value = 2757
$ws.Cells.Item(34, 6).Value = 1833
$ws.Cells.Item(36, 6).Value = 48
$ws.Cells.Item(37, 6).Value = 59
$ws.Cells.Item(38, 6).Value = 3690
$ws.Cells.Item(39, 6).Value = 319
$ws.Cells.Item(41, 6).Value = 16
$ws.Cells.Item(42, 6).Value = 44
$ws.Cells.Item(43, 6).Value = 918
$ws.Cells.Item(44, 6).Value = 539
$ws.Cells.Item(45, 6).Value = 34
$ws.Cells.Item(46, 6).Value = 1425
$ws.Cells.Item(47, 6).Value = 244
$ws.Cells.Item(50, 6).Value = 642
